$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.508.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.527.52"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -5.16%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.22"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -3.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.46"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.97%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.509"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.524.48"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.24%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.51%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.343"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.70%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.987.98"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -5.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "70.286.37"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000180"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.68%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "24.83"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -5.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.535.99"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.94%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.47"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -6.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.54"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -8.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "354.85"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -4.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.94"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -5.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.96"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.27%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.11"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.05"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -6.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.19"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -5.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.652.83"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.01"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0908"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.35%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "478.24"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.28"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.72%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.23%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.38"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.115"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +4.34%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.58"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.89%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.84"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.07%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.31"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.65"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -6.65%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.70"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.39"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.68%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "141.57"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -9.33%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.53"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -5.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.523"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -6.68%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -7.13%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.596"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.56%  "
